# "updated with new dashboards"
# Populates the mapping tables on the "D&C - Construction" and
# "Strategy & Operations" sheets (previously only header rows), and
# updates the active-sheet/selection view state accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: D&C - Construction  (xl/worksheets/sheet9.xml)
# ---------------------------------------------------------------------
$wsConstruction = $wb.Worksheets.Item("D&C - Construction")

# Row 1: C1/D1 keep their style but lose their shared-string value.
$wsConstruction.Range("C1").ClearContents()
$wsConstruction.Range("D1").ClearContents()

$constructionRows = @(
    @("AddDataPoint(D&C)", "Program"),
    @("AddDataPoint(D&C)", "Pending (Ha)"),
    @("DevCoAssessmentInput(D&C)", "Input Value"),
    @("AddDataPoint(D&C)", "Additional Data Point"),
    @("AddDataPoint(D&C)", "Value"),
    @("DevCoAssessmentAnalysis(D&C)", "Performance Signal Score"),
    @("DevCoAssessmentAnalysis(D&C)", "Assessment Criteria"),
    @("DevCoAssessmentAnalysis(D&C)", "Value"),
    @("DevCoAssessmentAnalysis(D&C)", "Rating")
)

$r = 2
foreach ($row in $constructionRows) {
    $wsConstruction.Cells.Item($r, 1).Value = $row[0]
    $wsConstruction.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Note: the runtime quantizes ColumnWidth to coarse pixel steps, so the
# inputs below are chosen to land as close as possible to the target
# stored widths of 29.44140625 / 12.109375 characters.
$wsConstruction.Columns.Item(1).ColumnWidth = 29.3
$wsConstruction.Columns.Item(2).ColumnWidth = 11.3

# ---------------------------------------------------------------------
# Sheet: Strategy & Operations  (xl/worksheets/sheet10.xml)
# ---------------------------------------------------------------------
$wsStrategy = $wb.Worksheets.Item("Strategy & Operations")

$wsStrategy.Range("C1").ClearContents()
$wsStrategy.Range("D1").ClearContents()

$strategyRows = @(
    @("DevCoAssessmentAnalysis(S&O)", "Rating"),
    @("DevCoAssessmentAnalysis(S&O)", "Assessment Criteria"),
    @("DevCoAssessmentAnalysis(S&O)", "Value"),
    @("AddDataPoint(S&O)", "Weightage"),
    @("AddDataPoint(S&O)", "Name of the KPI")
)

$r = 2
foreach ($row in $strategyRows) {
    $wsStrategy.Cells.Item($r, 1).Value = $row[0]
    $wsStrategy.Cells.Item($r, 2).Value = $row[1]
    $r++
}

$wsStrategy.Columns.Item(1).ColumnWidth = 29.3

# ---------------------------------------------------------------------
# View state: "Strategy & Operations" becomes the active/selected tab
# (moving away from "Innovation & Technology"), with a new selection on
# each touched sheet.
# ---------------------------------------------------------------------
$wsInnovation = $wb.Worksheets.Item("Innovation & Technology")
$wsInnovation.Activate()
$wsInnovation.Range("F23").Select()

$wsConstruction.Activate()
$wsConstruction.Range("A19").Select()

# Leave "Strategy & Operations" as the final active / selected tab.
$wsStrategy.Activate()
$wsStrategy.Range("C21").Select()

"done"
